$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.484.54"
$ws.Range("E2").Value = "  +5.71%  "
$ws.Range("D3").Value = "2.060.94"
$ws.Range("E3").Value = "  +4.49%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'252.78"
$ws.Range("E5").Value = "  +3.55%  "
$ws.Range("E6").Value = "  +2.87%  "
$ws.Range("D7").Value = "'66.16"
$ws.Range("E7").Value = "  +16.38%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.382"
$ws.Range("E9").Value = "  +6.92%  "
$ws.Range("D10").Value = "'59.51"
$ws.Range("E10").Value = "  +0.99%  "
$ws.Range("E11").Value = "  +5.52%  "
$ws.Range("E12").Value = "  +1.52%  "
$ws.Range("D13").Value = "'0.911"
$ws.Range("E13").Value = "  -2.68%  "
$ws.Range("D14").Value = "'14.97"
$ws.Range("E14").Value = "  +5.95%  "
$ws.Range("D15").Value = "2.363.33"
$ws.Range("E15").Value = "  +4.48%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "'21.23"
$ws.Range("E16").Value = "  +22.34%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "'5.58"
$ws.Range("E17").Value = "  +7.13%  "
$ws.Range("D18").Value = "2.064.83"
$ws.Range("E18").Value = "  +4.70%  "
$ws.Range("D19").Value = "37.280.15"
$ws.Range("E19").Value = "  +5.29%  "
$ws.Range("D20").Value = "'74.04"
$ws.Range("E20").Value = "  +4.40%  "
$ws.Range("D21").Value = "0.0₃0881"
$ws.Range("E21").Value = "  +4.99%  "
$ws.Range("D22").Value = "'5.49"
$ws.Range("E22").Value = "  +7.49%  "
$ws.Range("D23").Value = "'239.98"
$ws.Range("E23").Value = "  +3.42%  "
$ws.Range("E24").Value = "  +5.75%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "'2.39"
$ws.Range("E26").Value = "  +3.73%  "
$ws.Range("D27").Value = "'9.76"
$ws.Range("E27").Value = "  +8.01%  "
$ws.Range("E28").Value = "  -0.76%  "
$ws.Range("E29").Value = "  +4.30%  "
$ws.Range("D30").Value = "'5.28"
$ws.Range("E30").Value = "  +9.51%  "
$ws.Range("E31").Value = "  +3.48%  "
$ws.Range("D32").Value = "'0.114"
$ws.Range("E32").Value = "  +23.97%  "
$ws.Range("E33").Value = "  +6.85%  "
$ws.Range("D34").Value = "'4.77"
$ws.Range("E34").Value = "  +12.28%  "
$ws.Range("E35").Value = "  +5.23%  "
$ws.Range("D36").Value = "'2.46"
$ws.Range("E36").Value = "  +5.02%  "
$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'1.84"
$ws.Range("E38").Value = "  +4.63%  "
$ws.Range("B39").Value = "THORChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D39").Value = "'6.09"
$ws.Range("E39").Value = "  +19.65%  "
$ws.Range("D40").Value = "'3.08"
$ws.Range("E40").Value = "  +37.14%  "
$ws.Range("E41").Value = "  +17.03%  "
$ws.Range("E42").Value = "  +4.25%  "
$ws.Range("E43").Value = "  +5.08%  "
$ws.Range("E45").Value = "  +4.42%  "
$ws.Range("D46").Value = "'17.08"
$ws.Range("E46").Value = "  +7.53%  "
$ws.Range("D47").Value = "'95.51"
$ws.Range("E47").Value = "  +4.90%  "
$ws.Range("D48").Value = "'7.92"
$ws.Range("E48").Value = "  +5.87%  "
$ws.Range("D49").Value = "1.419.71"
$ws.Range("E49").Value = "  +3.17%  "
$ws.Range("E50").Value = "  +2.51%  "
$ws.Range("D51").Value = "'46.90"
